# DHW can now have its own heating system, independent of the room
# heating system. Insert a new "dhw heating system" column right after
# the existing "heating system" column (Z) - everything from the old
# "heat distribution" column onward (AA..AK) shifts one column to the
# right (AB..AL). The existing sample row's heating-system value
# ("GSHP") becomes the dhw heating system value, and the room heating
# system cell is reset to "None" (DHW-only example).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at AA (old column AA, "heat distribution",
# and everything after it, is pushed one column to the right).
$ws.Columns.Item(27).Insert()

# Header row (row 1) / unit-or-hint row (row 2) / sample-value row (row 3)
# for the newly inserted "dhw heating system" column.
$ws.Range("AA1").Value = "dhw heating system"
$ws.Range("AA2").Value = """same"" or the choices from heating system"
$ws.Range("AA3").Value = "GSHP"

# The room "heating system" sample value is now "None" (DHW supplied
# separately, no room heating system configured for this example row).
$ws.Range("Z3").Value = "None"

# Keep the selection on the newly inserted column for the sample row,
# matching where the edit was made.
$ws.Range("Z4").Select()
